$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "C3" = 329
    "D3" = 3162
    "F3" = 110.2
    "G3" = 1169.1
    "B4" = 66
    "C4" = 1317
    "D4" = 1345
    "E4" = 50.3
    "F4" = 447.4
    "G4" = 978.8
    "J4" = 1
    "M4" = 0.7
    "B5" = 194
    "C5" = 3483
    "D5" = 2478
    "E5" = 51.5
    "F5" = 540.8
    "G5" = 993.2
    "I5" = 2
    "L5" = 0.3
    "B6" = 447
    "C6" = 3684
    "D6" = 2495
    "E6" = 103
    "F6" = 665.6
    "G6" = 1392
    "I6" = 1
    "J6" = 30
    "L6" = 0.2
    "M6" = 16.7
    "B7" = 468
    "C7" = 2379
    "D7" = 1737
    "E7" = 87.5
    "F7" = 583
    "G7" = 1260.2
    "H7" = 3
    "I7" = 14
    "J7" = 69
    "K7" = 0.6
    "L7" = 3.4
    "M7" = 50.1
    "B8" = 274
    "C8" = 1015
    "D8" = 1017
    "E8" = 54.1
    "F8" = 398.1
    "G8" = 1079.8
    "H8" = 3
    "I8" = 17
    "J8" = 110
    "K8" = 0.6
    "L8" = 6.7
    "M8" = 116.8
    "B9" = 263
    "C9" = 638
    "D9" = 626
    "E9" = 49.6
    "F9" = 450.6
    "G9" = 964.7
    "I9" = 32
    "L9" = 22.6
    "M9" = 160.3
    "B10" = 199
    "C10" = 297
    "D10" = 351
    "E10" = 49.2
    "F10" = 409.9
    "G10" = 786.3
    "H10" = 13
    "I10" = 29
    "K10" = 3.2
    "L10" = 40
    "M10" = 179.2
    "B11" = 130
    "C11" = 208
    "D11" = 211
    "E11" = 74
    "F11" = 554.5
    "G11" = 978.1
    "H11" = 11
    "I11" = 23
    "J11" = 41
    "K11" = 6.3
    "L11" = 61.3
    "M11" = 190.1
    "B12" = 54
    "D12" = 67
    "E12" = 151.1
    "F12" = 453
    "G12" = 805.4
    "H12" = 1
    "I12" = 4
    "J12" = 13
    "K12" = 2.8
    "L12" = 31.2
    "M12" = 156.3
}

foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}
